$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell "Patio" in G1, matching the style of the existing F1 header cell
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Patio"
$excel.CutCopyMode = $false

# Update the active selection to match the committed state
$ws.Range("E19").Select()
